# Apply commit "feat: add 2022-Q4 data":
#   - 总计 (Total) sheet gains a new row for 2022-Q4 (inserted right after the header,
#     pushing the existing quarters down by one row).
#   - A brand-new "2022-Q4" worksheet is inserted right before "2022-Q3", holding the
#     quarterly fund-holdings breakdown for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: "总计" sheet - insert the 2022-Q4 summary row under the header.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Rows.Item(2).Insert()

# New row 2 contents.
$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q4"
$ws1.Range("C2").Value = 4
$ws1.Range("D2").Value = 0.1

# The inserted row picked up ad-hoc formatting; restore it to match the sheet's
# existing look: column A carries the bordered/centered "index" style, B:D plain.
$ws1.Range("A3").Copy()
$ws1.Range("A2").PasteSpecial(-4122)
$ws1.Range("B3:D3").Copy()
$ws1.Range("B2:D2").PasteSpecial(-4122)

# PasteSpecial only touched formatting, but re-assert values to be safe.
$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q4"
$ws1.Range("C2").Value = 4
$ws1.Range("D2").Value = 0.1

# Renumber the (0-based) index column for every row now that one more exists.
$ws1.Range("A3").Value = 1
$ws1.Range("A4").Value = 2
$ws1.Range("A5").Value = 3
$ws1.Range("A6").Value = 4

# ---------------------------------------------------------------------------
# Step 2: add the new "2022-Q4" worksheet right before "2022-Q3", with the same
# layout/styling as the other quarterly sheets.
# ---------------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($q3Sheet)
$newSheet.Name = "2022-Q4"

# Re-fetch "2022-Q3" by name: positional references re-bind after the insert.
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# Clone header row + the 4-row data block's styling from the 2022-Q3 sheet.
$q3Sheet.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$q3Sheet.Range("A2:H5").Copy($newSheet.Range("A2:H5"))

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'000965"
$newSheet.Range("C2").Value = "汇丰晋信新动力混合"
$newSheet.Range("D2").Value = "'0.95"
$newSheet.Range("E2").Value = "'91.04"
$newSheet.Range("F2").Value = "'5.55"
$newSheet.Range("G2").Value = "'0.0527"
$newSheet.Range("H2").Value = 2

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'011997"
$newSheet.Range("C3").Value = "景顺长城安盈回报一年持有期混合A"
$newSheet.Range("D3").Value = "'1.50"
$newSheet.Range("E3").Value = "'26.78"
$newSheet.Range("F3").Value = "'1.86"
$newSheet.Range("G3").Value = "'0.0279"
$newSheet.Range("H3").Value = 3

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'540004"
$newSheet.Range("C4").Value = "汇丰晋信2026周期混合"
$newSheet.Range("D4").Value = "'1.08"
$newSheet.Range("E4").Value = "'23.97"
$newSheet.Range("F4").Value = "'1.98"
$newSheet.Range("G4").Value = "'0.0214"
$newSheet.Range("H4").Value = 3

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'011998"
$newSheet.Range("C5").Value = "景顺长城安盈回报一年持有期混合C"
$newSheet.Range("D5").Value = "'0.08"
$newSheet.Range("E5").Value = "'26.78"
$newSheet.Range("F5").Value = "'1.86"
$newSheet.Range("G5").Value = "'0.0015"
$newSheet.Range("H5").Value = 3

# Keep the original active sheet/selection ("总计" was active before the edit).
$ws1.Activate()
$ws1.Range("A1").Select() | Out-Null

Write-Host "Inserted 2022-Q4 row in 总计 and added 2022-Q4 worksheet"
